$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (Price) and 1h volume change (Volume(1h)) cells.
# Values are written as text (apostrophe-prefixed Formula) and then restyled to
# the workbook default ("Normal") so no numeric auto-conversion or stray style
# index is introduced - matching the original inlineStr cells exactly.

$ws.Range("D2").Formula = "'37.868.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Formula = "'  -0.36%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Formula = "'2.030.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Formula = "'  -1.07%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Formula = "'  -0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Formula = "'228.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Formula = "'  -0.59%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Formula = "'  -0.60%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Formula = "'60.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Formula = "'  +2.69%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E9").Formula = "'  -1.12%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Formula = "'0.0815"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Formula = "'  +0.94%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Formula = "'  +0.38%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Formula = "'2.331.35"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Formula = "'  -1.06%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Formula = "'  -0.72%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Formula = "'21.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Formula = "'  +2.53%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Formula = "'0.761"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Formula = "'  +1.27%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Formula = "'  -2.30%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Formula = "'2.016.62"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Formula = "'  -1.61%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Formula = "'37.832.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Formula = "'  -0.29%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Formula = "'69.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Formula = "'  +0.23%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Formula = "'5.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Formula = "'  -5.18%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Formula = "'  -0.89%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Formula = "'224.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Formula = "'  -0.24%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Formula = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Formula = "'  -0.04%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E25").Formula = "'  +0.58%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Formula = "'166.95"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Formula = "'  +0.04%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Formula = "'  -4.38%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D30").Formula = "'1.28"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Formula = "'  -3.81%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Formula = "'  +0.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Formula = "'  +4.07%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Formula = "'4.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Formula = "'  -2.55%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Formula = "'  -0.74%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Formula = "'  -1.90%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Formula = "'6.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Formula = "'  +6.07%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Formula = "'2.28"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Formula = "'  -2.43%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Formula = "'3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Formula = "'  -1.00%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Formula = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Formula = "'  +0.01%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Formula = "'1.525.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Formula = "'  +2.54%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Formula = "'17.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Formula = "'  +3.86%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Formula = "'  +0.23%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Formula = "'96.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Formula = "'  -1.04%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Formula = "'  -1.03%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Formula = "'  -1.38%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Formula = "'  -1.66%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Formula = "'  -3.58%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Formula = "'  -0.77%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Formula = "'  -0.24%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Formula = "'  +0.79%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Formula = "'2.221.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Formula = "'  -1.02%  "
$ws.Range("E51").Style = "Normal"
